$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (column D) values, forcing text so Excel does not coerce these
#     dotted/numeric-looking strings into floating point numbers. ---
$dRows = @(2, 3, 5, 7, 8, 9, 12, 13, 14, 15, 16, 19, 21, 22, 23, 24, 25, 26, 27, 29, 30, 31, 32, 33, 34, 35, 36, 37, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50)
foreach ($r in $dRows) {
    $ws.Range("D$r").NumberFormat = "@"
}
$ws.Range("D2").Value = "25.764.10"
$ws.Range("D3").Value = "1.816.19"
$ws.Range("D5").Value = "277.57"
$ws.Range("D7").Value = "0.5075"
$ws.Range("D8").Value = "0.3531"
$ws.Range("D9").Value = "44.59"
$ws.Range("D12").Value = "0.8271"
$ws.Range("D13").Value = "0.07855"
$ws.Range("D14").Value = "1.825.30"
$ws.Range("D15").Value = "5.073"
$ws.Range("D16").Value = "87.67"
$ws.Range("D19").Value = "0.000008032"
$ws.Range("D21").Value = "25.804.06"
$ws.Range("D22").Value = "4.742"
$ws.Range("D23").Value = "10.01"
$ws.Range("D24").Value = "6.091"
$ws.Range("D25").Value = "142.58"
$ws.Range("D26").Value = "2.204"
$ws.Range("D27").Value = "1.675"
$ws.Range("D29").Value = "109.35"
$ws.Range("D30").Value = "4.344"
$ws.Range("D31").Value = "4.236"
$ws.Range("D32").Value = "0.08797"
$ws.Range("D33").Value = "0.04872"
$ws.Range("D34").Value = "0.7294"
$ws.Range("D35").Value = "1.138"
$ws.Range("D36").Value = "2.870"
$ws.Range("D37").Value = "1.0000"
$ws.Range("D39").Value = "2.387"
$ws.Range("D40").Value = "0.01850"
$ws.Range("D41").Value = "0.5165"
$ws.Range("D42").Value = "0.9643"
$ws.Range("D43").Value = "6.207"
$ws.Range("D44").Value = "110.61"
$ws.Range("D45").Value = "8.020"
$ws.Range("D46").Value = "1.000"
$ws.Range("D47").Value = "0.4570"
$ws.Range("D48").Value = "0.1366"
$ws.Range("D49").Value = "36.61"
$ws.Range("D50").Value = "9.213"
foreach ($r in $dRows) {
    $ws.Range("D$r").ClearFormats()
}

# --- Update Volume(1h) (column E) values ---
$ws.Range("E2").Value = "  -4.05%  "
$ws.Range("E3").Value = "  -3.07%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -7.86%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  -5.11%  "
$ws.Range("E8").Value = "  -5.91%  "
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("E10").Value = "  -7.18%  "
$ws.Range("E11").Value = "  -7.29%  "
$ws.Range("E12").Value = "  -6.95%  "
$ws.Range("E13").Value = "  -3.72%  "
$ws.Range("E14").Value = "  -2.79%  "
$ws.Range("E15").Value = "  -4.40%  "
$ws.Range("E16").Value = "  -6.28%  "
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("E18").Value = "  -4.71%  "
$ws.Range("E19").Value = "  -6.18%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("E21").Value = "  -4.15%  "
$ws.Range("E22").Value = "  -4.92%  "
$ws.Range("E23").Value = "  -6.10%  "
$ws.Range("E24").Value = "  -4.89%  "
$ws.Range("E25").Value = "  -2.64%  "
$ws.Range("E26").Value = "  -4.65%  "
$ws.Range("E27").Value = "  -3.30%  "
$ws.Range("E28").Value = "  -5.22%  "
$ws.Range("E30").Value = "  -8.23%  "
$ws.Range("E31").Value = "  -8.29%  "
$ws.Range("E32").Value = "  -3.81%  "
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("E34").Value = "  -10.99%  "
$ws.Range("E35").Value = "  -3.28%  "
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("E39").Value = "  -9.16%  "
$ws.Range("E40").Value = "  -5.15%  "
$ws.Range("E41").Value = "  -14.86%  "
$ws.Range("E42").Value = "  -9.92%  "
$ws.Range("E43").Value = "  -6.57%  "
$ws.Range("E44").Value = "  -3.99%  "
$ws.Range("E45").Value = "  -10.28%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  -10.90%  "
$ws.Range("E48").Value = "  -8.47%  "
$ws.Range("E49").Value = "  -2.62%  "
$ws.Range("E50").Value = "  -7.77%  "
$ws.Range("E51").Value = "  -7.73%  "

